$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = -138
$ws.Range("H2").Value = 1.01

$ws.Range("G3").Value = -16

$ws.Range("G4").Value = -70
$ws.Range("H4").Value = 1.05

$ws.Range("I5").Value = 0.28

$ws.Range("G6").Value = -70
$ws.Range("H6").Value = 1.05

$ws.Range("G7").Value = -430

$ws.Range("G8").Value = -430

$ws.Range("G9").Value = -430
